$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing comma formatting in column B (add space after comma, fix "budayamutu")
$ws.Range("B1").Value = "ekonomi, moneter"
$ws.Range("B2").Value = "ekonomi, mikro"
$ws.Range("B3").Value = "ekonomi"
$ws.Range("B4").Value = "corel draw"
$ws.Range("B5").Value = "kepemimpinan, organisasi"
$ws.Range("B6").Value = "manusia, salmon"
$ws.Range("B7").Value = "sekolah, jaringan"
$ws.Range("B8").Value = "kepemimpinan, sekolah, budaya mutu"
$ws.Range("B9").Value = "komunikasi, manusia"

# Append new rows 10-14 (repeat of rows 3-7)
$ws.Range("A10").Value = 3
$ws.Range("B10").Value = "ekonomi"

$ws.Range("A11").Value = 4
$ws.Range("B11").Value = "corel draw"

$ws.Range("A12").Value = 5
$ws.Range("B12").Value = "kepemimpinan, organisasi"

$ws.Range("A13").Value = 6
$ws.Range("B13").Value = "manusia, salmon"

$ws.Range("A14").Value = 7
$ws.Range("B14").Value = "sekolah, jaringan"

# Column C width (closest attainable to target 8.42578125 given engine's width quantization)
$ws.Columns.Item(3).ColumnWidth = 7.6

# Selection
$ws.Range("A10:B14").Select()
